# "Views with Encryption and CheckOption"
#
# Adds a new block of notes (rows 37-48) to Sheet1 documenting SQL Server
# views: view vs. table ownership/permission checks, SELECT access through
# a view vs. directly against the underlying/other tables, and the
# Account-based access outcomes.
#
# New shared strings must be created in the exact order they are introduced
# by the original edit so the rebuilt sharedStrings.xml lines up index-for-
# index with the target workbook (cells reference strings by index, and a
# brand-new string is appended to the shared-string table the first time its
# text is written). Order below mirrors that exactly - do not reorder these
# Value assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 37 / 38 - "view" block: who owns it, how it is qualified, how it's queried
$ws.Range("A37").Value = "view"
$ws.Range("D37").Value = "dbo.View"
$ws.Range("B38").Value = "dbo"
$ws.Range("C38").Value = "database owner"

# Row 41 / 42 - "table" block: two underlying tables, same denial outcome
$ws.Range("A41").Value = "table"
$ws.Range("C41").Value = "dbo.Table1"
$ws.Range("C42").Value = "dbo.Table2"

$ws.Range("B37").Value = "dbo.view"
$ws.Range("C37").Value = "select"

$ws.Range("B45").Value = "select * from dbo.View"

$ws.Range("D41").Value = "denied select"

$ws.Range("B47").Value = "Have I got dbo.View?"
$ws.Range("B48").Value = "Within dbo.View, I can access to ALL dbo objects"

$ws.Range("F41").Value = "select * from dbo.Table1"

# Row 43 / 44 - Account-based access comparison
$ws.Range("C43").Value = "Account.Table3"
$ws.Range("D43").Value = " if denied select"
$ws.Range("F43").Value = "no access with view"
$ws.Range("D44").Value = "if ok by Account"
$ws.Range("F44").Value = "access ok"

# Cells 42/B, 42/D, 42/F repeat values already introduced above (reuse the
# same shared-string entries - no new unique strings from here on).
$ws.Range("B42").Value = "table"
$ws.Range("D42").Value = "denied select"
$ws.Range("F42").Value = "select * from dbo.Table1"

# Match the saved file's final selection/active cell.
$ws.Range("F44").Select()
